$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '36.734.72'
Set-TextValue $ws.Range('E2') '  +1.86%  '
Set-TextValue $ws.Range('D3') '1.957.04'
Set-TextValue $ws.Range('E3') '  +1.97%  '
Set-TextValue $ws.Range('E4') '  -0.13%  '
Set-TextValue $ws.Range('D5') '243.83'
Set-TextValue $ws.Range('E5') '  +1.42%  '
Set-TextValue $ws.Range('D6') '0.614'
Set-TextValue $ws.Range('E6') '  +2.07%  '
Set-TextValue $ws.Range('D7') '58.33'
Set-TextValue $ws.Range('E7') '  +5.59%  '
Set-TextValue $ws.Range('E8') '  -0.11%  '
Set-TextValue $ws.Range('D9') '0.371'
Set-TextValue $ws.Range('E9') '  +3.08%  '
Set-TextValue $ws.Range('D10') '0.0807'
Set-TextValue $ws.Range('E10') '  -1.40%  '
Set-TextValue $ws.Range('E11') '  +0.38%  '
Set-TextValue $ws.Range('D12') '22.18'
Set-TextValue $ws.Range('E12') '  +7.50%  '
Set-TextValue $ws.Range('D13') '2.244.71'
Set-TextValue $ws.Range('E13') '  +1.85%  '
Set-TextValue $ws.Range('D14') '0.819'
Set-TextValue $ws.Range('E14') '  +2.07%  '
Set-TextValue $ws.Range('D15') '13.63'
Set-TextValue $ws.Range('E15') '  +4.03%  '
Set-TextValue $ws.Range('D16') '5.27'
Set-TextValue $ws.Range('E16') '  +2.58%  '
Set-TextValue $ws.Range('D17') '1.960.96'
Set-TextValue $ws.Range('E17') '  +1.08%  '
Set-TextValue $ws.Range('D18') '36.708.05'
Set-TextValue $ws.Range('E18') '  +1.91%  '
Set-TextValue $ws.Range('D19') '69.69'
Set-TextValue $ws.Range('E19') '  +1.96%  '
Set-TextValue $ws.Range('E20') '  +1.00%  '
Set-TextValue $ws.Range('D21') '5.11'
Set-TextValue $ws.Range('E21') '  +4.35%  '
Set-TextValue $ws.Range('D22') '228.16'
Set-TextValue $ws.Range('E22') '  +1.23%  '
Set-TextValue $ws.Range('E23') '  -0.02%  '
Set-TextValue $ws.Range('D24') '2.40'
Set-TextValue $ws.Range('E24') '  -0.48%  '
Set-TextValue $ws.Range('E25') '  +3.98%  '
Set-TextValue $ws.Range('D26') '9.30'
Set-TextValue $ws.Range('E26') '  +0.86%  '
Set-TextValue $ws.Range('B27') 'Kaspa'
Set-TextValue $ws.Range('C27') 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D27') '0.138'
Set-TextValue $ws.Range('E27') '  +17.07%  '
Set-TextValue $ws.Range('B28') 'Monero'
Set-TextValue $ws.Range('C28') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D28') '160.78'
Set-TextValue $ws.Range('E28') '  -0.49%  '
Set-TextValue $ws.Range('D29') '19.34'
Set-TextValue $ws.Range('E29') '  +1.60%  '
Set-TextValue $ws.Range('E30') '  +2.18%  '
Set-TextValue $ws.Range('D31') '1.11'
Set-TextValue $ws.Range('E31') '  +0.06%  '
Set-TextValue $ws.Range('D32') '4.67'
Set-TextValue $ws.Range('E32') '  +1.84%  '
Set-TextValue $ws.Range('D33') '0.0619'
Set-TextValue $ws.Range('E33') '  +0.36%  '
Set-TextValue $ws.Range('E34') '  +0.31%  '
Set-TextValue $ws.Range('D35') '6.27'
Set-TextValue $ws.Range('E35') '  +6.44%  '
Set-TextValue $ws.Range('B36') 'BinanceUSD'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range('D36') '1.00'
Set-TextValue $ws.Range('E36') '  -0.24%  '
Set-TextValue $ws.Range('B37') 'RenderToken'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D37') '3.43'
Set-TextValue $ws.Range('E37') '  +20.59%  '
Set-TextValue $ws.Range('D38') '2.19'
Set-TextValue $ws.Range('E38') '  +4.17%  '
Set-TextValue $ws.Range('E39') '  -0.54%  '
Set-TextValue $ws.Range('D40') '0.0997'
Set-TextValue $ws.Range('E40') '  +5.26%  '
Set-TextValue $ws.Range('D41') '2.90'
Set-TextValue $ws.Range('E41') '  +2.67%  '
Set-TextValue $ws.Range('E42') '  +3.44%  '
Set-TextValue $ws.Range('E43') '  +0.87%  '
Set-TextValue $ws.Range('D44') '16.07'
Set-TextValue $ws.Range('E44') '  +4.32%  '
Set-TextValue $ws.Range('E45') '  +2.06%  '
Set-TextValue $ws.Range('D46') '1.345.47'
Set-TextValue $ws.Range('E46') '  +1.43%  '
Set-TextValue $ws.Range('D47') '87.46'
Set-TextValue $ws.Range('E47') '  +1.18%  '
Set-TextValue $ws.Range('D48') '7.15'
Set-TextValue $ws.Range('E48') '  +0.32%  '
Set-TextValue $ws.Range('E49') '  +1.37%  '
Set-TextValue $ws.Range('D50') '2.137.92'
Set-TextValue $ws.Range('E50') '  +1.92%  '
Set-TextValue $ws.Range('B51') 'MultiversX'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue $ws.Range('D51') '43.46'
Set-TextValue $ws.Range('E51') '  -2.89%  '
